$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.3897772753026277
$ws.Range("J2").Value = 0.3897772753026277
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 7.082880333333333
$ws.Range("N2").Value = 21.248641
$ws.Range("O2").Value = 0.4363153076427043
$ws.Range("P2").Value = 0.4363153076427044
$ws.Range("Q2").Value = 2.845636890400888
$ws.Range("R2").Value = 25.610732013608
$ws.Range("S2").Value = 0.1700657917858011
$ws.Range("T2").Value = 0.1700657917858011

# Row 3
$ws.Range("I3").Value = 0.3897772753026277
$ws.Range("J3").Value = 0.3897772753026277
$ws.Range("O3").Value = 0.04522408903652051
$ws.Range("P3").Value = 0.04522408903652051
$ws.Range("S3").Value = 0.0176273222026984
$ws.Range("T3").Value = 0.0176273222026984

# Row 4
$ws.Range("I4").Value = 0.3897772753026277
$ws.Range("J4").Value = 0.3897772753026277
$ws.Range("M4").Value = 7.805874333333333
$ws.Range("N4").Value = 23.417623
$ws.Range("O4").Value = 0.4808527464653325
$ws.Range("P4").Value = 0.4808527464653325
$ws.Range("Q4").Value = 3.136108887824888
$ws.Range("R4").Value = 28.224979990424
$ws.Range("S4").Value = 0.1874254733390426
$ws.Range("T4").Value = 0.1874254733390426

# Row 5
$ws.Range("I5").Value = 0.3897772753026277
$ws.Range("J5").Value = 0.3897772753026277
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.6105033333333334
$ws.Range("N5").Value = 1.83151
$ws.Range("O5").Value = 0.03760785685544264
$ws.Range("P5").Value = 0.03760785685544264
$ws.Range("Q5").Value = 0.2452774472088889
$ws.Range("R5").Value = 2.20749702488
$ws.Range("S5").Value = 0.01465868797508568
$ws.Range("T5").Value = 0.01465868797508568

# Row 6
$ws.Range("G6").Value = 0.6289866666666667
$ws.Range("H6").Value = 1.88696
$ws.Range("I6").Value = 0.6102227246973723
$ws.Range("J6").Value = 0.6102227246973724
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 7.082880333333333
$ws.Range("N6").Value = 21.248641
$ws.Range("O6").Value = 0.4363153076427043
$ws.Range("P6").Value = 0.4363153076427044
$ws.Range("Q6").Value = 4.455037291262222
$ws.Range("R6").Value = 40.09533562136
$ws.Range("S6").Value = 0.2662495158569033
$ws.Range("T6").Value = 0.2662495158569033

# Row 7
$ws.Range("G7").Value = 0.6289866666666667
$ws.Range("H7").Value = 1.88696
$ws.Range("I7").Value = 0.6102227246973723
$ws.Range("J7").Value = 0.6102227246973724
$ws.Range("O7").Value = 0.04522408903652051
$ws.Range("P7").Value = 0.04522408903652051
$ws.Range("Q7").Value = 0.4617646907911112
$ws.Range("R7").Value = 4.155882217119999
$ws.Range("S7").Value = 0.02759676683382211
$ws.Range("T7").Value = 0.02759676683382211

# Row 8
$ws.Range("G8").Value = 0.6289866666666667
$ws.Range("H8").Value = 1.88696
$ws.Range("I8").Value = 0.6102227246973723
$ws.Range("J8").Value = 0.6102227246973724
$ws.Range("M8").Value = 7.805874333333333
$ws.Range("N8").Value = 23.417623
$ws.Range("O8").Value = 0.4808527464653325
$ws.Range("P8").Value = 0.4808527464653325
$ws.Range("Q8").Value = 4.909790877342222
$ws.Range("R8").Value = 44.18811789607999
$ws.Range("S8").Value = 0.2934272731262899
$ws.Range("T8").Value = 0.29342727312629

# Row 9
$ws.Range("G9").Value = 0.6289866666666667
$ws.Range("H9").Value = 1.88696
$ws.Range("I9").Value = 0.6102227246973723
$ws.Range("J9").Value = 0.6102227246973724
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.6105033333333334
$ws.Range("N9").Value = 1.83151
$ws.Range("O9").Value = 0.03760785685544264
$ws.Range("P9").Value = 0.03760785685544264
$ws.Range("Q9").Value = 0.3839984566222223
$ws.Range("R9").Value = 3.4559861096
$ws.Range("S9").Value = 0.02294916888035696
$ws.Range("T9").Value = 0.02294916888035696

